# Trade #50 closed at 2026-02-17 12:49:03 - unknown UNKNOWN +0.000%
#
# This script applies three related updates to the live trading results
# workbook:
#   1. Summary sheet       - bump Total Trades (B6) and recompute Win Rate % (B9)
#   2. Strategy Status sheet - bump MarketMaking Trades (D4) and Win Rate % (G4)
#   3. All Trades / MarketMaking sheets - append the new trade as row 51

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Summary sheet
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B6").Value = 50      # Total Trades
$summary.Range("B9").Value = 42      # Win Rate %

# ---------------------------------------------------------------------
# 2. Strategy Status sheet (MarketMaking row)
# ---------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("D4").Value = 50       # Trades
$status.Range("G4").Value = 42       # Win Rate %

# ---------------------------------------------------------------------
# 3. Append the new trade row (row 51) to both "All Trades" and
#    "MarketMaking" sheets, which mirror each other.
# ---------------------------------------------------------------------
$sheetNames = @("All Trades", "MarketMaking")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Copy the Date cell (B2) into B51: every row in this column holds
    # the literal text "2026-02-17" (same value the new trade needs), and
    # Copy/Paste preserves the text type without Excel re-interpreting it
    # as a date serial value (which a direct .Value assignment would do).
    $ws.Range("B2").Copy($ws.Range("B51"))

    $ws.Range("A51").Value = 50
    $ws.Range("C51").Value = "12:48:57"
    $ws.Range("D51").Value = "MarketMaking"
    $ws.Range("E51").Value = "UP"
    $ws.Range("F51").Value = 0.96
    $ws.Range("G51").Value = 0.96
    $ws.Range("H51").Value = "CLOSED"
    $ws.Range("I51").Value = 0
    $ws.Range("J51").Value = 0
    $ws.Range("K51").Value = 100.16
    $ws.Range("L51").Value = 0
    $ws.Range("M51").Value = 0
    $ws.Range("N51").Value = 0.6
    $ws.Range("O51").Value = "Normal spread capture: 19600 bps"
    $ws.Range("P51").Value = "early_exit"
    $ws.Range("Q51").Value = 0.13
}
